# Auto-generated Excel COM-interop edit script
# Applies value updates to the Asura_Profits-style market data sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the target diff.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 529.8
$ws.Range("J6").Value = 999.5
$ws.Range("L6").Value = 2998.5
$ws.Range("N6").Value = -3222.5
$ws.Range("H8").Value = 6714.5557
$ws.Range("I8").Value = 61.57143
$ws.Range("J8").Value = 30000
$ws.Range("K8").Value = 184.71429
$ws.Range("L8").Value = 90000
$ws.Range("M8").Value = -45.71429000000001
$ws.Range("N8").Value = -90278
$ws.Range("H9").Value = 158.33333
$ws.Range("I9").Value = 158.33333
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 158.33333
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 10.66667000000001
$ws.Range("N9").ClearContents()
$ws.Range("H64").Value = 3326.5757
$ws.Range("I64").Value = 3232.25
$ws.Range("J64").Value = 3380.476
$ws.Range("K64").Value = 3232.25
$ws.Range("L64").Value = 3380.476
$ws.Range("M64").Value = -2984.25
$ws.Range("N64").Value = -3876.476
$ws.Range("H67").Value = 3326.5757
$ws.Range("I67").Value = 3232.25
$ws.Range("J67").Value = 3380.476
$ws.Range("K67").Value = 3232.25
$ws.Range("L67").Value = 3380.476
$ws.Range("M67").Value = -2374.25
$ws.Range("N67").Value = -5096.476000000001
$ws.Range("H76").Value = 4528.5713
$ws.Range("I76").Value = 4611.1113
$ws.Range("J76").Value = 4380
$ws.Range("K76").Value = 4611.1113
$ws.Range("L76").Value = 4380
$ws.Range("M76").Value = -4296.1113
$ws.Range("N76").Value = -5010
$ws.Range("H79").Value = 4528.5713
$ws.Range("I79").Value = 4611.1113
$ws.Range("J79").Value = 4380
$ws.Range("K79").Value = 4611.1113
$ws.Range("L79").Value = 4380
$ws.Range("M79").Value = -3519.1113
$ws.Range("N79").Value = -6564
$ws.Range("H112").Value = 2173.8076
$ws.Range("J112").Value = 2348.6956
$ws.Range("L112").Value = 7046.0868
$ws.Range("N112").Value = -9262.086800000001
$ws.Range("H133").Value = 69808.336
$ws.Range("J133").Value = 69808.336
$ws.Range("L133").Value = 69808.336
$ws.Range("N133").Value = -79928.336
$ws.Range("H134").Value = 125160
$ws.Range("J134").Value = 125160
$ws.Range("L134").Value = 125160
$ws.Range("N134").Value = -135300
$ws.Range("H137").Value = 1402.4464
$ws.Range("I137").Value = 1404.1333
$ws.Range("J137").Value = 1400.5
$ws.Range("K137").Value = 4212.3999
$ws.Range("L137").Value = 4201.5
$ws.Range("M137").Value = -1662.3999
$ws.Range("N137").Value = -9301.5
$ws.Range("H138").Value = 2084.7
$ws.Range("I138").Value = 1345.9412
$ws.Range("J138").Value = 2236.012
$ws.Range("K138").Value = 4037.8236
$ws.Range("L138").Value = 6708.036
$ws.Range("M138").Value = 1102.1764
$ws.Range("N138").Value = -16988.036

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1169.1578
$ws.Range("I2").Value = 750.0714
$ws.Range("K2").Value = 750.0714
$ws.Range("M2").Value = -637.0714
$ws.Range("H63").Value = 3132
$ws.Range("I63").Value = 2804.9524
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 2804.9524
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = -2118.9524
$ws.Range("N63").Value = -11372
$ws.Range("H66").Value = 3132
$ws.Range("I66").Value = 2804.9524
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 14024.762
$ws.Range("L66").Value = 50000
$ws.Range("M66").Value = -10592.762
$ws.Range("N66").Value = -56864
$ws.Range("H116").Value = 1169.1578
$ws.Range("I116").Value = 750.0714
$ws.Range("K116").Value = 750.0714
$ws.Range("M116").Value = 1543.9286
$ws.Range("H122").Value = 2347.5454
$ws.Range("I122").Value = 2347.5454
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7042.6362
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4592.6362
$ws.Range("N122").ClearContents()
$ws.Range("H133").Value = 39025.777
$ws.Range("J133").Value = 39025.777
$ws.Range("L133").Value = 39025.777
$ws.Range("N133").Value = -44085.777
$ws.Range("H134").Value = 69228
$ws.Range("J134").Value = 69228
$ws.Range("L134").Value = 69228
$ws.Range("N134").Value = -79368

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1169.1578
$ws.Range("I3").Value = 750.0714
$ws.Range("K3").Value = 750.0714
$ws.Range("M3").Value = -636.0714
$ws.Range("H107").Value = 28576.1
$ws.Range("I107").Value = 32960.117
$ws.Range("J107").Value = 3733.3333
$ws.Range("K107").Value = 32960.117
$ws.Range("L107").Value = 3733.3333
$ws.Range("M107").Value = -31040.117
$ws.Range("N107").Value = -7573.3333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 8999.538
$ws.Range("J4").Value = 8999.538
$ws.Range("L4").Value = 8999.538
$ws.Range("N4").Value = -9223.538
$ws.Range("H5").Value = 166669000
$ws.Range("I5").Value = 500000260
$ws.Range("J5").Value = 3360
$ws.Range("K5").Value = 500000260
$ws.Range("L5").Value = 3360
$ws.Range("M5").Value = -500000148
$ws.Range("N5").Value = -3584
$ws.Range("H62").Value = 60123.332
$ws.Range("I62").Value = 85701.664
$ws.Range("J62").Value = 8966.666999999999
$ws.Range("K62").Value = 85701.664
$ws.Range("L62").Value = 8966.666999999999
$ws.Range("M62").Value = -85077.664
$ws.Range("N62").Value = -10214.667
$ws.Range("H65").Value = 60123.332
$ws.Range("I65").Value = 85701.664
$ws.Range("J65").Value = 8966.666999999999
$ws.Range("K65").Value = 428508.32
$ws.Range("L65").Value = 44833.335
$ws.Range("M65").Value = -425388.32
$ws.Range("N65").Value = -51073.335

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 899
$ws.Range("I8").Value = 899
$ws.Range("K8").Value = 2697
$ws.Range("M8").Value = -2558
$ws.Range("H68").Value = 1235.4166
$ws.Range("J68").Value = 1203
$ws.Range("L68").Value = 3609
$ws.Range("N68").Value = -5231
$ws.Range("H71").Value = 1235.4166
$ws.Range("J71").Value = 1203
$ws.Range("L71").Value = 10827
$ws.Range("N71").Value = -18939
$ws.Range("H137").Value = 47620600
$ws.Range("I137").Value = 765
$ws.Range("J137").Value = 111113710
$ws.Range("K137").Value = 2295
$ws.Range("L137").Value = 333341130
$ws.Range("M137").Value = 2805
$ws.Range("N137").Value = -333351330

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").ClearContents()
$ws.Range("H38").Value = 30000
$ws.Range("H40").Value = 50000
$ws.Range("J40").Value = 50000
$ws.Range("L40").Value = 50000
$ws.Range("N40").Value = -50302
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H80").Value = 3062.5
$ws.Range("J80").Value = 3166.6667
$ws.Range("L80").Value = 3166.6667
$ws.Range("N80").Value = -5162.6667
$ws.Range("H83").Value = 3062.5
$ws.Range("J83").Value = 3166.6667
$ws.Range("L83").Value = 15833.3335
$ws.Range("N83").Value = -25817.3335
$ws.Range("H126").Value = 2072.9473
$ws.Range("I126").Value = 1881.6471
$ws.Range("J126").Value = 3699
$ws.Range("K126").Value = 5644.9413
$ws.Range("L126").Value = 11097
$ws.Range("M126").Value = -3174.9413
$ws.Range("N126").Value = -16037
$ws.Range("H140").Value = 48898.184
$ws.Range("J140").Value = 48898.184
$ws.Range("L140").Value = 48898.184
$ws.Range("N140").Value = -59258.184

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9630.75
$ws.Range("I40").Value = 11841
$ws.Range("K40").Value = 11841
$ws.Range("M40").Value = -11705
$ws.Range("H46").Value = 1800.1333
$ws.Range("I46").Value = 1566.6666
$ws.Range("J46").Value = 1955.7778
$ws.Range("K46").Value = 1566.6666
$ws.Range("L46").Value = 1955.7778
$ws.Range("M46").Value = -1378.6666
$ws.Range("N46").Value = -2331.7778
$ws.Range("H82").Value = 2181.1875
$ws.Range("I82").Value = 1971.8
$ws.Range("J82").Value = 2530.1667
$ws.Range("K82").Value = 1971.8
$ws.Range("L82").Value = 2530.1667
$ws.Range("M82").Value = -1610.8
$ws.Range("N82").Value = -3252.1667
$ws.Range("H85").Value = 2181.1875
$ws.Range("I85").Value = 1971.8
$ws.Range("J85").Value = 2530.1667
$ws.Range("K85").Value = 1971.8
$ws.Range("L85").Value = 2530.1667
$ws.Range("M85").Value = -723.8
$ws.Range("N85").Value = -5026.1667
$ws.Range("H122").Value = 10874088
$ws.Range("I122").Value = 16671706
$ws.Range("J122").Value = 3556.25
$ws.Range("K122").Value = 50015118
$ws.Range("L122").Value = 10668.75
$ws.Range("M122").Value = -50012668
$ws.Range("N122").Value = -15568.75
$ws.Range("H132").Value = 4973.7646
$ws.Range("I132").Value = 4842.839
$ws.Range("K132").Value = 14528.517
$ws.Range("M132").Value = -11998.517
$ws.Range("H134").Value = 78919.664
$ws.Range("J134").Value = 78919.664
$ws.Range("L134").Value = 78919.664
$ws.Range("N134").Value = -89059.664

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 27779428
$ws.Range("I122").Value = 50001400
$ws.Range("J122").Value = 1962.5
$ws.Range("K122").Value = 150004200
$ws.Range("L122").Value = 5887.5
$ws.Range("M122").Value = -150001750
$ws.Range("N122").Value = -10787.5
$ws.Range("H126").Value = 5124.684
$ws.Range("I126").Value = 8497.1
$ws.Range("J126").Value = 1377.5555
$ws.Range("K126").Value = 25491.3
$ws.Range("L126").Value = 4132.666499999999
$ws.Range("M126").Value = -23021.3
$ws.Range("N126").Value = -9072.666499999999
$ws.Range("H132").Value = 2448.2942
$ws.Range("I132").Value = 2916.6155
$ws.Range("J132").Value = 926.25
$ws.Range("K132").Value = 8749.8465
$ws.Range("L132").Value = 2778.75
$ws.Range("M132").Value = -6219.8465
$ws.Range("N132").Value = -7838.75
